# benerin sweetallert di delete koordinat sama template excel
# -> fix the "delete coordinate" sweetalert + the Excel template.
#
# The template sheet is renamed to the generic "Sheet1", and the
# Latitude/Longitude/Tarif sample columns are converted from raw numeric
# placeholders into the literal text values that the app now writes
# (coordinates as formatted text like "-6.358211", and the tariff column
# holding a status code like "P3" instead of a bare rupiah number).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet used to be named after a to-do note; rename it to the normal
# default sheet name used by the template going forward.
$ws.Name = "Sheet1"

# Latitude/Longitude sample rows (2 & 3): these were stored as raw numbers
# (e.g. -6358211 / 106815949, missing their decimal points). Re-enter them
# as proper decimal text so they round-trip as literal strings instead of
# being re-interpreted as numbers. Formatting the range as Text first keeps
# Excel from re-parsing the numeric-looking strings back into numbers.
$ws.Range("A2:B3").NumberFormat = "@"
$ws.Range("B2").Value = "106.815949"
$ws.Range("B3").Value = "106.813574"
$ws.Range("A2").Value = "-6.358211"
$ws.Range("A3").Value = "-6.375125"

# Tarif column (I2/I3): replace the old numeric rupiah placeholders with the
# "P3" status code text used by the current template.
$ws.Range("I2").Value = "P3"
$ws.Range("I3").Value = "P3"

# The trailing "..." placeholder row (row 4, columns A & B) moves onto the
# same text-formatted style as the coordinate columns above.
$ws.Range("A4:B4").NumberFormat = "@"

# Match the template's saved selection.
$ws.Range("B6").Select()
